$d = $word.ActiveDocument
$p = $d.Paragraphs.Add()
$r = $p.Range
$f = $d.Fields.Add($r, -1, "MERGEFIELD =metadata.generator \* MERGEFORMAT", $false)
Write-Output "field count: $($d.Fields.Count)"
$f2 = $d.Fields.Item($d.Fields.Count)
$res = $f2.Result
Write-Output "result: $($res.Start) $($res.End) [$($res.Text)]"
$res.NoProofing = $true
Write-Output "set noproofing"
